# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51)
# with refreshed figures, and swaps the Bittensor/Hedera rows (48/49).
#
# Each entry is: A1 cell reference, new text value, and whether the value
# must be pinned to Text format (otherwise Excel's COM layer would infer a
# plain numeric string such as "0.999" as a real number instead of text).
$updates = @(
    @('D2', '53.627.43', $false),
    @('E2', '  -4.88%  ', $false),
    @('D3', '2.223.63', $false),
    @('E3', '  -6.47%  ', $false),
    @('D4', '0.999', $true),
    @('E4', '  -0.16%  ', $false),
    @('D5', '482.37', $true),
    @('E5', '  -3.98%  ', $false),
    @('D6', '125.65', $true),
    @('E6', '  -3.66%  ', $false),
    @('D7', '0.999', $true),
    @('E7', '  +0.02%  ', $false),
    @('D8', '0.517', $true),
    @('D9', '2.228.65', $false),
    @('E9', '  -6.54%  ', $false),
    @('D10', '0.0916', $true),
    @('E10', '  -6.89%  ', $false),
    @('E11', '  -1.57%  ', $false),
    @('D12', '4.68', $true),
    @('E12', '  -1.63%  ', $false),
    @('D13', '0.315', $true),
    @('E13', '  -3.11%  ', $false),
    @('D14', '2.618.81', $false),
    @('E14', '  -6.57%  ', $false),
    @('D15', '21.05', $true),
    @('E15', '  -2.37%  ', $false),
    @('D16', '53.545.03', $false),
    @('E16', '  -5.00%  ', $false),
    @('E17', '  -3.89%  ', $false),
    @('D18', '2.212.07', $false),
    @('E18', '  -8.09%  ', $false),
    @('D19', '9.56', $true),
    @('E19', '  -4.82%  ', $false),
    @('D21', '297.67', $true),
    @('E21', '  -3.24%  ', $false),
    @('E22', '  -2.76%  ', $false),
    @('E23', '  +0.30%  ', $false),
    @('D24', '63.27', $true),
    @('E24', '  -3.37%  ', $false),
    @('E25', '  -0.01%  ', $false),
    @('D26', '0.362', $true),
    @('E26', '  -1.52%  ', $false),
    @('D27', '0.142', $true),
    @('E27', '  -4.13%  ', $false),
    @('D28', '6.96', $true),
    @('E28', '  -4.20%  ', $false),
    @('D29', '169.37', $true),
    @('E29', '  -1.37%  ', $false),
    @('D30', '0.0₃0678', $false),
    @('E30', '  -5.24%  ', $false),
    @('E31', '  -4.24%  ', $false),
    @('E32', '  -0.11%  ', $false),
    @('D33', '0.997', $true),
    @('E33', '  -0.06%  ', $false),
    @('D34', '5.71', $true),
    @('E34', '  -0.95%  ', $false),
    @('D35', '1.04', $true),
    @('E35', '  -4.04%  ', $false),
    @('D36', '17.36', $true),
    @('E36', '  -1.32%  ', $false),
    @('D37', '1.14', $true),
    @('E37', '  -2.68%  ', $false),
    @('D38', '0.832', $true),
    @('E38', '  +4.43%  ', $false),
    @('E39', '  -5.51%  ', $false),
    @('D40', '35.66', $true),
    @('E40', '  -1.26%  ', $false),
    @('E41', '  -1.23%  ', $false),
    @('E42', '  -1.94%  ', $false),
    @('E43', '  -2.67%  ', $false),
    @('D44', '122.69', $true),
    @('E44', '  -6.30%  ', $false),
    @('E45', '  -2.65%  ', $false),
    @('E46', '  -3.60%  ', $false),
    @('E47', '  -5.68%  ', $false),
    @('B48', 'Hedera', $false),
    @('C48', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', $false),
    @('D48', '0.0468', $true),
    @('E48', '  -2.80%  ', $false),
    @('B49', 'Bittensor', $false),
    @('C49', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', $false),
    @('D49', '228.81', $true),
    @('E49', '  -5.17%  ', $false),
    @('E50', '  -3.37%  ', $false),
    @('D51', '15.97', $true),
    @('E51', '  -5.70%  ', $false)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $forceText = $update[2]

    $cell = $ws.Range($cellRef)
    if ($forceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $newValue
}
